# Updates cryptos list: price (D) and 1h volume change (E) columns,
# plus a swap of the Polygon / WrappedEther rows (12 & 13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.121.46"
$ws.Range("E2").Value = "  +0.95%  "

$ws.Range("D3").Value = "1.890.55"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7389"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.29%  "

$ws.Range("E8").Value = "  +1.74%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.85"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.99%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07168"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08345"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.62%  "

$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7583"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.27%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.905.25"
$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.422"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.07"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.156"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.36%  "

$ws.Range("D17").Value = "30.122.36"
$ws.Range("E17").Value = "  +0.99%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "251.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007863"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.72%  "

$ws.Range("D21").Value = "2.177.51"
$ws.Range("E21").Value = "  +2.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.946"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.72%  "

$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("E25").Value = "  -1.85%  "

$ws.Range("E26").Value = "  -0.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.056"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.479"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.581"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.45%  "

$ws.Range("E32").Value = "  +0.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.201"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05359"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.256"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7713"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.730"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01961"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.71%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.762"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4570"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.42%  "

$ws.Range("D42").Value = "1.100.73"
$ws.Range("E42").Value = "  +0.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.072"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.52%  "

$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8748"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.004"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.861"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.589"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.601"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.62%  "

$ws.Range("D51").Value = "2.067.00"
$ws.Range("E51").Value = "  +1.46%  "
